# Insert a new data row at row 34 (pushing the existing rows 34..111 down
# to 35..112) and populate it with the new weekly price record.
#
# This mirrors the diff: every row from 34 onward in the "before" file
# reappears one row lower in the "after" file, and the dimension grows
# from A1:T111 to A1:T112. The newly opened row 34 receives the values
# below (date 2021-10-29 / serial 44498, Región Metropolitana, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 34:111 down to 35:112, leaving a blank row 34 behind.
$ws.Rows.Item(34).Insert()

# Fill the newly inserted row 34 with the new record's data.
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44498
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100101
$ws.Range("H34").Value = "Berries"
$ws.Range("I34").Value = 100101001
$ws.Range("J34").Value = "Arándano (blue)"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 380
$ws.Range("N34").Value = 12000
$ws.Range("O34").Value = 12000
$ws.Range("P34").Value = 12000
$ws.Range("Q34").Value = "$/bandeja 2 kilos"
$ws.Range("R34").Value = "Región Metropolitana"
$ws.Range("S34").Value = 6000
$ws.Range("T34").Value = 2
